# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 04:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1095210
$ws.Range("C4").Value = 187
$ws.Range("D4").Value = 155737
$ws.Range("E4").Value = 875612
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 63861

# Row 27 - Pakistan
$ws.Range("B27").Value = 16817
$ws.Range("C27").Value = 344
$ws.Range("E27").Value = 12327
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = 385

# Row 38 - Corea del Sur
$ws.Range("B38").Value = 10774
$ws.Range("C38").Value = 9
$ws.Range("D38").Value = 9072
$ws.Range("E38").Value = 1454
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 248
